$wb = $excel.ActiveWorkbook

# --- Sheet1: add a "Seniority Level" column (G) with data ------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("G1").Value = "Seniority Level"
$ws1.Range("G2").Value = "Mid-Level"
$ws1.Range("G3").Value = "Entry-Level"
$ws1.Range("G4").Value = "Senior-Level"

# Row 2 picked up a slightly taller custom row height.
$ws1.Rows.Item(2).RowHeight = 17

# New column G needs a bit more width (stored <col width="19".../> — the
# COM ColumnWidth property is ~0.8333 narrower than the stored XML width).
$ws1.Columns.Item(7).ColumnWidth = 18.16666666666667

# Put a thin grid border around the whole (now wider) table.
$ws1.Range("A1:G4").Borders.LineStyle = 1

# --- Sheet2: no data change, only an incidental style reindex happens ------
# (handled automatically by the style table when Sheet1's formatting is
# written above; nothing else to do for Sheet2.)

# --- New blank worksheet "Sheet5", placed right after Sheet2 ---------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws5 = $wb.Worksheets.Add($null, $ws2)
$ws5.Name = "Sheet5"
$ws5.Range("I31:I32").Select() | Out-Null

# --- Restore Sheet1 as the active/selected sheet, with its new view state --
$ws1.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 200
$ws1.Range("E9").Select() | Out-Null
